$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.346.59'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '3.107.47'
$ws.Range("E3").Value = '  -2.64%  '
$ws.Range("E4").Value = '  +0.71%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '589.16'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.32%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '151.90'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.45%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.104.62'
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("E10").Value = '  -1.37%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '5.93'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.87%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.459'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("E13").Value = '  +2.09%  '
$ws.Range("E14").Value = '  -2.22%  '
$ws.Range("D15").Value = '3.628.50'
$ws.Range("E15").Value = '  -2.36%  '
$ws.Range("E16").Value = '  -1.86%  '
$ws.Range("E17").Value = '  +1.96%  '
$ws.Range("D18").Value = '63.918.92'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").Value = '3.109.50'
$ws.Range("E19").Value = '  -1.51%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '469.12'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.67%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '14.87'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +3.78%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.737'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.41%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '7.57'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.72%  '
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("E25").Value = '  +5.81%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '81.72'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("E27").Value = '  -0.31%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '9.72'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +3.87%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '7.40'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +3.78%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.70'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.53%  '
$ws.Range("E31").Value = '  +0.70%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.67%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.116'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +6.00%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '27.45'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.68%  '
$ws.Range("D35").Value = '0.0₃0846'
$ws.Range("E35").Value = '  -1.83%  '
$ws.Range("E36").Value = '  +0.83%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '6.15'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +2.08%  '
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.35'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("E39").Value = '  -2.86%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '9.35'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +5.17%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '50.83'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.11%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '451.40'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +2.71%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.290'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.58%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0370'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("D45").Value = '2.847.99'
$ws.Range("E45").Value = '  -2.46%  '
$ws.Range("E46").Value = '  +0.90%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '38.35'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.94%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '129.90'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.40%  '
$ws.Range("E49").Value = '  +0.03%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '25.15'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +3.86%  '
$ws.Range("E51").Value = '  +3.73%  '
